# Generate Report for Handback
# Mirrors the localization hand-back run: the status text flips from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet
# that shows it, the per-locale "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns get populated now that the
# hand-back has happened, and a couple of columns get widened so the new
# (longer) values aren't clipped.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdName    = "82a0a388-0d8e-455e-9b26-ad91f2b58d3c.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/812df44cdc10d13b2cfd11095ab98148ddae40f6/e2e/82a0a388-0d8e-455e-9b26-ad91f2b58d3c.md"
$hyperlinkColor = 15570276  # OLE BGR for FF6495ED (cornflower blue), matching the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: the per-locale status columns mirror the same string,
# so they flip to the new status too.
# ---------------------------------------------------------------------
$overview.Cells.Item(2, 5).Value = $newStatus   # E2 (zh-cn status)
$overview.Cells.Item(2, 6).Value = $newStatus   # F2 (de-de status)

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Cells.Item(2, 3).Value = $newStatus   # C2 Status

# I2 Latest Target File -> hyperlink to the source .md, same display text
# used by the A2 hyperlink.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $hyperlinkColor

# J2 Latest Handback File -> the zh-cn xlf that was handed back
$zhcn.Cells.Item(2, 10).Value = "82a0a388-0d8e-455e-9b26-ad91f2b58d3c.4c8ae8774d7a21fb69634e20ae7bd7a7108515a5.zh-cn.xlf"

# K2 Latest Handback DateTime
$zhcn.Cells.Item(2, 11).Value = "2016-09-03 15:03:51"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Cells.Item(2, 3).Value = $newStatus   # C2 Status

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $hyperlinkColor

# J2 Latest Handback File -> the de-de xlf that was handed back
$dede.Cells.Item(2, 10).Value = "82a0a388-0d8e-455e-9b26-ad91f2b58d3c.4c8ae8774d7a21fb69634e20ae7bd7a7108515a5.de-de.xlf"

# K2 Latest Handback DateTime (de-de finished a few seconds after zh-cn)
$dede.Cells.Item(2, 11).Value = "2016-09-03 15:03:58"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15
